$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '50.723.63'
$ws.Range("E2").Value = '  -1.80%  '

$ws.Range("D3").Value = '2.918.71'
$ws.Range("E3").Value = '  -2.04%  '

$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '369.46'
$ws.Range("E5").Value = '  -2.80%  '

$ws.Range("D6").Value = '99.57'
$ws.Range("E6").Value = '  -5.49%  '

$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  -2.33%  '

$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").Value = '0.571'
$ws.Range("E9").Value = '  -3.73%  '

$ws.Range("D10").Value = '35.38'
$ws.Range("E10").Value = '  -5.31%  '

$ws.Range("D11").Value = '0.138'
$ws.Range("E11").Value = '  -1.12%  '

$ws.Range("D12").Value = '0.0832'
$ws.Range("E12").Value = '  -1.54%  '

$ws.Range("D13").Value = '3.380.14'
$ws.Range("E13").Value = '  -2.10%  '

$ws.Range("D14").Value = '17.69'
$ws.Range("E14").Value = '  -4.05%  '

$ws.Range("D15").Value = '7.31'
$ws.Range("E15").Value = '  -3.35%  '

$ws.Range("D16").Value = '11.25'
$ws.Range("E16").Value = '  +51.64%  '

$ws.Range("D17").Value = '2.923.48'
$ws.Range("E17").Value = '  -1.93%  '

$ws.Range("D18").Value = '0.948'
$ws.Range("E18").Value = '  -2.57%  '

$ws.Range("D19").Value = '50.683.21'
$ws.Range("E19").Value = '  -1.84%  '

$ws.Range("D20").Value = '3.09'
$ws.Range("E20").Value = '  -7.65%  '

$ws.Range("D21").Value = '12.14'
$ws.Range("E21").Value = '  -6.61%  '

$ws.Range("D22").Value = '0.0₃0940'
$ws.Range("E22").Value = '  -2.53%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").Value = '67.94'
$ws.Range("E23").Value = '  -2.33%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D24").Value = '260.98'
$ws.Range("E24").Value = '  -0.45%  '

$ws.Range("D25").Value = '3.02'
$ws.Range("E25").Value = '  +6.49%  '

$ws.Range("D26").Value = '8.08'
$ws.Range("E26").Value = '  +2.74%  '

$ws.Range("D27").Value = '7.19'
$ws.Range("E27").Value = '  -5.71%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("D29").Value = '0.164'
$ws.Range("E29").Value = '  -5.24%  '

$ws.Range("D30").Value = '25.19'
$ws.Range("E30").Value = '  -2.87%  '

$ws.Range("D31").Value = '0.108'
$ws.Range("E31").Value = '  -3.67%  '

$ws.Range("D32").Value = '9.77'
$ws.Range("E32").Value = '  -1.54%  '

$ws.Range("D33").Value = '50.31'
$ws.Range("E33").Value = '  -1.21%  '

$ws.Range("D34").Value = '2.02'
$ws.Range("E34").Value = '  -3.01%  '

$ws.Range("D35").Value = '0.0433'
$ws.Range("E35").Value = '  -3.52%  '

$ws.Range("D36").Value = '32.15'
$ws.Range("E36").Value = '  -8.98%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").Value = '3.10'
$ws.Range("E38").Value = '  +1.25%  '

$ws.Range("D39").Value = '0.113'
$ws.Range("E39").Value = '  -2.21%  '

$ws.Range("D40").Value = '15.96'
$ws.Range("E40").Value = '  -7.43%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.75'
$ws.Range("E41").Value = '  -5.71%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '2.42'
$ws.Range("E42").Value = '  -8.30%  '

$ws.Range("D43").Value = '118.87'
$ws.Range("E43").Value = '  -5.06%  '

$ws.Range("D44").Value = '20.65'
$ws.Range("E44").Value = '  -6.00%  '

$ws.Range("D45").Value = '2.02'
$ws.Range("E45").Value = '  -1.37%  '

$ws.Range("D46").Value = '0.267'
$ws.Range("E46").Value = '  -8.52%  '

$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '2.28'
$ws.Range("E47").Value = '  -3.72%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '3.19'
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("D49").Value = '1.963.69'
$ws.Range("E49").Value = '  -3.91%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '3.231.37'
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").Value = '0.0320'
$ws.Range("E51").Value = '  -6.54%  '

